# Daily update of covid19 tracker data files
# Bump the "date" column (column B) from 43933 (2020-04-12) to
# 43934 (2020-04-13) for every data row (5 through 96) on the
# "Country Updates" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Country Updates")

$ws.Range("B5:B96").Value = 43934
